$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "datum" -> "date": rename the header of column B.
$ws.Range("B1").Value = "date"

# Clarify the B1 ("date") comment with an example.
$ws.Range("B1").CommentThreaded.Text("Tag und Uhrzeit des Zeitabschnitts. Beispiel: `"Fr 08:00`" oder `"Sa 10:15`".")

# Add a new threaded comment explaining the "priority" column (D1).
$ws.Range("D1").AddCommentThreaded("Priorität, nach der ein Zeitabschnitt nicht verwendet werden soll. Höhere Zahl = Timeslot wird eher freigehalten.")

# Restore the active selection.
$ws.Range("G16").Select()
